# Rik Activity 2019.xlsx — "Baguette eating...updated data through 3/31"
# Appends 5 new activity-log rows (149-153) to the "2019" sheet / Table2,
# covering 3/31-4/1 2019, and updates the view selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019")
$lo = $ws.ListObjects.Item("Table2")

# Formula used by Table2's calculated "Z" column.
$zFormula = '=IF(Table2[[#This Row],[Activity]]="Sleep",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,"NA")'

# Helper: append one data row to the bottom of Table2, growing the table
# range/autofilter via ListRows.Add(), then writing the real cell data
# (ListRows.Add() alone does not materialize sheet cells in this runtime).
# NOTE: positional params only — named param binding (-Start ...) does not
# reliably pass values through in this PS host, so call this positionally.
function Add-ActivityRow {
    param(
        [double]$Start,
        $End,
        [string]$Activity,
        $Comment
    )

    $lo.ListRows.Add() | Out-Null
    $r = $lo.Range.Rows.Count + $lo.Range.Row - 1

    # Column A (Start) - set value, then copy formatting from the row above.
    $ws.Cells.Item($r, 1).Value = $Start
    $ws.Cells.Item($r - 1, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null

    if ($null -ne $End) {
        $ws.Cells.Item($r, 2).Value = $End
        $ws.Cells.Item($r - 1, 1).Copy()
        $ws.Cells.Item($r, 2).PasteSpecial(-4122) | Out-Null
    }

    $ws.Cells.Item($r, 3).Value = $Activity

    if ($null -ne $Comment) {
        $ws.Cells.Item($r, 4).Value = $Comment
    }

    $ws.Cells.Item($r, 5).Formula = $zFormula
    $ws.Cells.Item($r - 1, 5).Copy()
    $ws.Cells.Item($r, 5).PasteSpecial(-4122) | Out-Null

    return $r
}

Add-ActivityRow 43555.456250000003 $null "Food" "Baguette" | Out-Null
Add-ActivityRow 43555.506944444445 $null "Food" "Hamburger" | Out-Null
Add-ActivityRow 43555.791666666664 $null "Food" "Chicken, coconut rice" | Out-Null
Add-ActivityRow 43555.913263888891 43556.236111111109 "Sleep" $null | Out-Null
$lastRow = Add-ActivityRow 43556.270833333336 $null "Food" "Latte"

$excel.ActiveWindow.ScrollRow = 118
$ws.Range("E" + $lastRow).Select() | Out-Null
